$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "307.41"
Set-TextValue $ws.Range("E2") "6.98%"

Set-TextValue $ws.Range("D3") "31.82"
Set-TextValue $ws.Range("E3") "7.53%"

Set-TextValue $ws.Range("D4") "5.297"
Set-TextValue $ws.Range("E4") "3.26%"

Set-TextValue $ws.Range("D5") "0.07336"
Set-TextValue $ws.Range("E5") "9.43%"

Set-TextValue $ws.Range("D6") "7.851"
Set-TextValue $ws.Range("E6") "7.03%"

Set-TextValue $ws.Range("D7") "3.771"
Set-TextValue $ws.Range("E7") "11.05%"

Set-TextValue $ws.Range("D8") "1.489"
Set-TextValue $ws.Range("E8") "8.90%"

Set-TextValue $ws.Range("D9") "0.9060"
Set-TextValue $ws.Range("E9") "-1.57%"

Set-TextValue $ws.Range("D10") "0.01625"
Set-TextValue $ws.Range("E10") "2,418.46%"

Set-TextValue $ws.Range("D11") "0.1688"
Set-TextValue $ws.Range("E11") "6.04%"

Set-TextValue $ws.Range("D12") "0.07516"
Set-TextValue $ws.Range("E12") "10.05%"

Set-TextValue $ws.Range("D13") "0.08049"
Set-TextValue $ws.Range("E13") "4.22%"

Set-TextValue $ws.Range("D14") "0.03036"
Set-TextValue $ws.Range("E14") "3.34%"

Set-TextValue $ws.Range("D15") "0.09975"
Set-TextValue $ws.Range("E15") "11.01%"

Set-TextValue $ws.Range("D16") "0.001515"
Set-TextValue $ws.Range("E16") "-4.65%"

Set-TextValue $ws.Range("D17") "0.04578"
Set-TextValue $ws.Range("E17") "2.15%"

Set-TextValue $ws.Range("D18") "0.006309"
Set-TextValue $ws.Range("E18") "0.92%"

Set-TextValue $ws.Range("D19") "3.481"
Set-TextValue $ws.Range("E19") "1.03%"

Set-TextValue $ws.Range("D20") "2.232"
Set-TextValue $ws.Range("E20") "0.19%"

Set-TextValue $ws.Range("D21") "0.3326"
Set-TextValue $ws.Range("E21") "3.47%"

Set-TextValue $ws.Range("D22") "0.1349"
Set-TextValue $ws.Range("E22") "3.55%"

Set-TextValue $ws.Range("D23") "4.326"
Set-TextValue $ws.Range("E23") "6.48%"

Set-TextValue $ws.Range("D24") "0.1644"
Set-TextValue $ws.Range("E24") "3.98%"

Set-TextValue $ws.Range("D25") "0.001229"
Set-TextValue $ws.Range("E25") "3.09%"

Set-TextValue $ws.Range("D26") "0.004434"
Set-TextValue $ws.Range("E26") "7.58%"

Set-TextValue $ws.Range("D27") "0.0001308"
Set-TextValue $ws.Range("E27") "9.11%"

Set-TextValue $ws.Range("D28") "0.0001759"
Set-TextValue $ws.Range("E28") "8.87%"

Set-TextValue $ws.Range("D40") "0.04506"
Set-TextValue $ws.Range("E40") "5.55%"

Set-TextValue $ws.Range("D41") "0.007243"
Set-TextValue $ws.Range("E41") "7.56%"

Set-TextValue $ws.Range("D42") "0.1347"
Set-TextValue $ws.Range("E42") "8.57%"

Set-TextValue $ws.Range("D43") "0.002274"
Set-TextValue $ws.Range("E43") "3.46%"

Set-TextValue $ws.Range("D44") "0.01411"
Set-TextValue $ws.Range("E44") "17.70%"

Set-TextValue $ws.Range("D45") "0.00006090"
Set-TextValue $ws.Range("E45") "7.14%"

Set-TextValue $ws.Range("E46") "-3.82%"

Set-TextValue $ws.Range("D47") "0.01314"
Set-TextValue $ws.Range("E47") "0.64%"
